$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-19.18267561737553, 1.978161539253551, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553),
    @(-19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, 1.546604066755419, -19.18267561737553, -19.18267561737553),
    @(-19.18267561737553, 2.000158186889835, 1.687894549519026, -19.18267561737553, 3.396654362913008, -19.18267561737553, -19.18267561737553, -19.18267561737553, 0.7666452887858605, -19.18267561737553),
    @(-19.18267561737553, 1.677639226893133, -19.18267561737553, -19.18267561737553, -19.18267561737553, 2.646528757336752, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553),
    @(-19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553),
    @(2.457248400598794, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553),
    @(-19.18267561737553, -19.18267561737553, -19.18267561737553, 1.794727997296051, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553),
    @(3.858805240872973, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553),
    @(-19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, 1.773326822249918, -19.18267561737553, 2.287928151373216),
    @(-19.18267561737553, -19.18267561737553, -19.18267561737553, 2.910561714938326, -19.18267561737553, 2.898581874806351, -19.18267561737553, -19.18267561737553, -19.18267561737553, 1.949937696921655),
    @(-19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553),
    @(-19.18267561737553, -19.18267561737553, -19.18267561737553, 2.518013427604846, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, 1.673495317358611, 1.844156449378137),
    @(-19.18267561737553, -19.18267561737553, 1.532978932485157, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, 1.759630547399855),
    @(-19.18267561737553, -19.18267561737553, 1.783686897404215, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553),
    @(-19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, 1.953322304095891, -19.18267561737553),
    @(-19.18267561737553, 2.133285979837471, 1.796585582879678, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, 2.006205961145961, 2.488948805052731, -19.18267561737553),
    @(-19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, 4.321925791655112, 2.033419522165497, 2.491201471714662, -19.18267561737553),
    @(-19.18267561737553, -19.18267561737553, 2.04506311327043, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553, 1.644180835999274, -19.18267561737553, -19.18267561737553),
    @(-19.18267561737553, 1.036508268327281, 1.506070703674997, -19.18267561737553, 3.24311422706878, -19.18267561737553, -19.18267561737553, 1.277145192404896, -19.18267561737553, 2.096366425069013),
    @(-19.18267561737553, 1.276304648189434, -19.18267561737553, 1.715206020153174, -19.18267561737553, 2.651098153083538, -19.18267561737553, -19.18267561737553, -19.18267561737553, -19.18267561737553)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowvals = $data[$r]
    for ($c = 0; $c -lt $rowvals.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $rowvals[$c]
    }
}